# Add new columns I ("I0") and J ("IF") to Sheet1, mirroring the
# existing header/data layout (A:H), and populate their values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Copy the formatting from the existing header cell H1 (bold, bordered,
# centered) onto the two new header cells so they match the rest of the
# header row, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (2-17) ---
$iVals = @(8, 4, 5, 6, 9, 5, 4, 4, 5, 7, 7, 6, 6, 5, 6, 5)
$jVals = @(8, 4, 5, 7, 9, 6, 4, 4, 6, 7, 7, 6, 6, 6, 6, 5)

for ($idx = 0; $idx -lt 16; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}

$excel.CutCopyMode = 0
